$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.261.78"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "3.561.12"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.33"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.18"
$ws.Range("E6").Value = "  +3.81%  "
$ws.Range("D7").Value = "3.560.09"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("E10").Value = "  +5.84%  "
$ws.Range("E11").Value = "  +7.65%  "
$ws.Range("E12").Value = "  +3.73%  "
$ws.Range("E13").Value = "  +2.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.26"
$ws.Range("E14").Value = "  +5.47%  "
$ws.Range("D15").Value = "4.163.36"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "3.562.67"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").Value = "68.259.49"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("E19").Value = "  +5.65%  "
$ws.Range("E20").Value = "  +6.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.07"
$ws.Range("E21").Value = "  +12.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "454.09"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("E23").Value = "  +4.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.47"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("D26").Value = "3.704.22"
$ws.Range("E26").Value = "  +1.79%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  +13.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.50"
$ws.Range("E29").Value = "  +4.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.72"
$ws.Range("E30").Value = "  +11.66%  "
$ws.Range("E31").Value = "  +3.68%  "
$ws.Range("E32").Value = "  +4.06%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.40"
$ws.Range("E34").Value = "  +5.46%  "
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("E36").Value = "  +4.91%  "
$ws.Range("D37").Value = "3.554.73"
$ws.Range("E37").Value = "  +1.86%  "
$ws.Range("E38").Value = "  +3.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.39"
$ws.Range("E39").Value = "  +8.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "180.89"
$ws.Range("E41").Value = "  +3.78%  "
$ws.Range("E42").Value = "  +5.01%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +4.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.29"
$ws.Range("E45").Value = "  +14.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.897"
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.23"
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("E48").Value = "  +5.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.66"
$ws.Range("E49").Value = "  +4.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.77"
$ws.Range("E50").Value = "  +3.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.263"
$ws.Range("E51").Value = "  +7.65%  "
